{"js": "// Append two new paragraphs after the existing content:\n//   1. \"Dkjs,,cn cx\" -- with \"cn\" flagged by the proofing engine\n//      (wrapped in <w:proofErr w:type=\"spellStart\"/>...<w:proofErr w:type=\"spellEnd\"/>)\n//   2. \",mm ccx\" followed by a manual line break (<w:br/>)\n//\n// body.insertParagraph() is used first (instead of insertOoxml directly) so\n// Word creates a proper new <w:p> block; inserting fresh empty paragraphs\n// also avoids inheriting run formatting (e.g. the superscript \"th\" run\n// earlier in the document) the way a plain insertParagraph(text, ...) call\n// would. We then overwrite each new paragraph's contents with an exact\n// WordprocessingML fragment (via insertOoxml/Replace) so the proofErr\n// markers and the line break land precisely as in the target markup.\n\nconst body = context.document.body;\n\nfunction flatOpcParagraph(innerXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' + innerXml + '</w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\nconst para1Xml = flatOpcParagraph(\n  '<w:r><w:t>Dkjs,,</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>cn</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> cx</w:t></w:r>'\n);\n\nconst para2Xml = flatOpcParagraph(\n  '<w:r><w:t>,mm ccx</w:t></w:r>' +\n  '<w:r><w:br/></w:r>'\n);\n\nconst p1 = body.insertParagraph(\"\", Word.InsertLocation.end);\np1.insertOoxml(para1Xml, Word.InsertLocation.replace);\n\nconst p2 = body.insertParagraph(\"\", Word.InsertLocation.end);\np2.insertOoxml(para2Xml, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Collapse to the very end of the document body (after the last paragraph's\n# text, before the final section break) and inject the two new paragraphs as\n# raw WordprocessingML so we get full control over run/proofErr/break\n# structure without inheriting formatting from the preceding run.\n$rng = $d.Content\n$rng.Collapse(0)\n\n$xml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n         '<w:r><w:t>Dkjs,,</w:t></w:r>' +\n         '<w:proofErr w:type=\"spellStart\"/>' +\n         '<w:r><w:t>cn</w:t></w:r>' +\n         '<w:proofErr w:type=\"spellEnd\"/>' +\n         '<w:r><w:t xml:space=\"preserve\"> cx</w:t></w:r>' +\n       '</w:p>' +\n       '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n         '<w:r><w:t>,mm ccx</w:t></w:r>' +\n         '<w:r><w:br/></w:r>' +\n       '</w:p>'\n\n$rng.InsertXML($xml)\n"}
